$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2022-01-17", "overview", "K02000001", "United Kingdom", 15305410, 84429, 85, 152075),
    @("2022-01-18", "overview", "K02000001", "United Kingdom", 15399300, 94432, 438, 152513),
    @("2022-01-19", "overview", "K02000001", "United Kingdom", 15506750, 108069, 359, 152872),
    @("2022-01-20", "overview", "K02000001", "United Kingdom", 15613283, 107364, 330, 153202),
    @("2022-01-21", "overview", "K02000001", "United Kingdom", 15709059, 95787, 288, 153490),
    @("2022-01-22", "overview", "K02000001", "United Kingdom", 15784488, 76807, 297, 153787)
)

$startRow = 524
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).Style = "Normal"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
